$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new "pair_kind" column (J) for the practice rows 2-5 ---
# These rows already had kind "generic" values elsewhere (col C); the sheet
# gains a parallel pair_kind entry of "generic" in column J.
$ws.Range("J2:J5").Value = "generic"

# --- New "stim details" block appended below the existing table ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# Rows 29-32: "video" word_type entries for months 6,6,7,7
$videoMonths = @(6, 6, 7, 7)
for ($i = 0; $i -lt $videoMonths.Length; $i++) {
    $r = 29 + $i
    $ws.Cells.Item($r, 1).Value = $videoMonths[$i]
    $ws.Cells.Item($r, 2).Value = "video"
}

# Rows 33-36: "audio" word_type entries for months 6,6,7,7
$audioMonths = @(6, 6, 7, 7)
for ($i = 0; $i -lt $audioMonths.Length; $i++) {
    $r = 33 + $i
    $ws.Cells.Item($r, 1).Value = $audioMonths[$i]
    $ws.Cells.Item($r, 2).Value = "audio"
}
